$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price (D) column whose new value could be misread as a plain
# number by Excel (e.g. '1.00' -> 1, '5.40' -> 5.4), losing the original
# textual formatting. Force each such cell to Text format individually
# before writing the value (a unioned/multi-area Range does not reliably
# propagate NumberFormat in this runtime).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '62.769.79'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '2.460.20'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '573.41'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '2.460.13'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').Value = '28.88'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '2.905.86'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = '62.753.33'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '2.464.35'
$ws.Range('D19').Value = '7.98'
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').Value = '10.99'
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').Value = '326.70'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +10.37%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '10.16'
$ws.Range('E25').Value = '  +20.12%  '
$ws.Range('D26').Value = '65.57'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').Value = '658.47'
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.584.85'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0975'
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -14.67%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').Value = '7.98'
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('D33').Value = '1.84'
$ws.Range('E33').Value = '  -1.72%  '
$ws.Range('E34').Value = '  -2.94%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '1.53'
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('D37').Value = '4.74'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '0.368'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '5.40'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E44').Value = '  -67.67%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '152.61'
$ws.Range('E46').Value = '  +5.15%  '
$ws.Range('D47').Value = '15.22'
$ws.Range('E47').Value = '  +1.72%  '
$ws.Range('D48').Value = '3.58'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').Value = '20.52'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '0.0511'
$ws.Range('E51').Value = '  -0.73%  '
